$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store plain text that often LOOKS like a
# number (e.g. "240.28", "  -0.03%  "). A bare .Value assignment would let
# Excel auto-coerce those into real numbers/percentages, which would change
# the stored cell type from a string to a number and diverge from the
# original inline-string layout. Force the ranges to Text format first so the
# values are kept verbatim as strings, then restore the "Normal" style so we
# don't leave a stray number-format applied to the cells.
$priceCol = $ws.Range("D2:D51")
$volCol = $ws.Range("E2:E51")
$priceCol.NumberFormat = "@"
$volCol.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.382.37"

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.848.74"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "240.28"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.30%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.07629"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.11%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +0.77%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07745"
$ws.Range("E11").Value = "  +0.04%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +0.43%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "0.6784"
$ws.Range("E13").Value = "  -0.16%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "0.00001057"
$ws.Range("E14").Value = "  -2.97%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "83.21"

# Row 16 - Uniswap
$ws.Range("E16").Value = "  +0.05%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.424.82"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "227.93"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "12.35"
$ws.Range("E19").Value = "  -0.78%  "

# Row 20 - Dai
$ws.Range("D20").Value = "0.9998"

# Row 21 - Chainlink
$ws.Range("D21").Value = "7.510"
$ws.Range("E21").Value = "  +0.60%  "

# Row 22 - BinanceUSD
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 - Monero
$ws.Range("D23").Value = "158.50"
$ws.Range("E23").Value = "  +0.61%  "

# Row 24 - Stellar
$ws.Range("E24").Value = "  -0.57%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "8.401"
$ws.Range("E25").Value = "  +0.34%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "17.71"
$ws.Range("E26").Value = "  +0.32%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "1.383"
$ws.Range("E27").Value = "  +5.90%  "

# Row 28 - PancakeSwap
$ws.Range("D28").Value = "1.462"
$ws.Range("E28").Value = "  -0.35%  "

# Row 29 - Hedera
$ws.Range("D29").Value = "0.05591"
$ws.Range("E29").Value = "  -0.70%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "4.109"
$ws.Range("E30").Value = "  -0.18%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "4.078"
$ws.Range("E31").Value = "  +0.93%  "

# Row 32 - LidoDAOToken
$ws.Range("D32").Value = "1.838"
$ws.Range("E32").Value = "  -0.79%  "

# Row 33 - ARBITRUM (unchanged)

# Row 34 - ImmutableX
$ws.Range("D34").Value = "0.6952"
$ws.Range("E34").Value = "  -2.04%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.582"
$ws.Range("E35").Value = "  -0.12%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.01805"
$ws.Range("E36").Value = "  +0.42%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.231.99"
$ws.Range("E37").Value = "  -0.30%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "2.713"
$ws.Range("E38").Value = "  -2.27%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "6.385"
$ws.Range("E39").Value = "  -1.39%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "0.9059"
$ws.Range("E40").Value = "  -0.30%  "

# Row 41 - PaxDollar
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - Quant
$ws.Range("D42").Value = "101.50"
$ws.Range("E42").Value = "  +0.09%  "

# Row 43 - Aave
$ws.Range("D43").Value = "66.03"
$ws.Range("E43").Value = "  -0.06%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "7.189"
$ws.Range("E44").Value = "  +0.42%  "

# Row 45 - TheSandbox
$ws.Range("D45").Value = "0.4008"
$ws.Range("E45").Value = "  -0.17%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "9.037"
$ws.Range("E46").Value = "  -0.07%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "1.676"
$ws.Range("E47").Value = "  -0.63%  "

# Row 48 - Algorand
$ws.Range("D48").Value = "0.1134"
$ws.Range("E48").Value = "  +0.93%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "0.05701"
$ws.Range("E49").Value = "  -0.15%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  +0.09%  "

# Row 51 - SynthetixNetwork -> NEARProtocol
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.346"
$ws.Range("E51").Value = "  +0.28%  "

# Restore the default "Normal" style on the touched columns so no stray
# number-format index is left referenced on the cells themselves.
$priceCol.Style = "Normal"
$volCol.Style = "Normal"
